$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Periodo Mora" column (E16:E23) - previous account statement periods are
# removed and replaced with the new set of periods, in reverse order.
$periodos = @("2005", "2004", "2003", "2002", "2001", "1912", "1911", "1910")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}

# "Valor Mora" column (F16:F23) - row 16 and row 23 values are swapped.
$ws.Range("F16").Value = 21533
$ws.Range("F17").Value = 34000
$ws.Range("F18").Value = 34000
$ws.Range("F19").Value = 34000
$ws.Range("F20").Value = 34000
$ws.Range("F21").Value = 34000
$ws.Range("F22").Value = 34000
$ws.Range("F23").Value = 34000

# "Salario Basico" column (G16:G23) - database value updated for all rows.
$ws.Range("G16:G23").Value = 877803
